$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Angptl3"
$ws.Cells.Item(2,3).Value = "Itga5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.030331
$ws.Cells.Item(2,8).Value = 0.090993
$ws.Cells.Item(2,9).Value = 0.005348220947889333
$ws.Cells.Item(2,10).Value = 0.005348220947889334
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 45.90594266666667
$ws.Cells.Item(2,14).Value = 137.717828
$ws.Cells.Item(2,15).Value = 0.3954672001633582
$ws.Cells.Item(2,16).Value = 0.3954672001633583
$ws.Cells.Item(2,17).Value = 1.392373147022667
$ws.Cells.Item(2,18).Value = 12.531358323204
$ws.Cells.Item(2,19).Value = 0.002115045964116816
$ws.Cells.Item(2,20).Value = 0.002115045964116817
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Angptl3"
$ws.Cells.Item(3,3).Value = "Itga5"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.030331
$ws.Cells.Item(3,8).Value = 0.090993
$ws.Cells.Item(3,9).Value = 0.005348220947889333
$ws.Cells.Item(3,10).Value = 0.005348220947889334
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 40.44578266666667
$ws.Cells.Item(3,14).Value = 121.337348
$ws.Cells.Item(3,15).Value = 0.3484294080560655
$ws.Cells.Item(3,16).Value = 0.3484294080560656
$ws.Cells.Item(3,17).Value = 1.226761034062667
$ws.Cells.Item(3,18).Value = 11.040849306564
$ws.Cells.Item(3,19).Value = 0.00186347745902613
$ws.Cells.Item(3,20).Value = 0.001863477459026131
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Angptl3"
$ws.Cells.Item(4,3).Value = "Itga5"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.030331
$ws.Cells.Item(4,8).Value = 0.090993
$ws.Cells.Item(4,9).Value = 0.005348220947889333
$ws.Cells.Item(4,10).Value = 0.005348220947889334
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 12.761795
$ws.Cells.Item(4,14).Value = 38.28538500000001
$ws.Cells.Item(4,15).Value = 0.1099393900775594
$ws.Cells.Item(4,16).Value = 0.1099393900775594
$ws.Cells.Item(4,17).Value = 0.387078004145
$ws.Cells.Item(4,18).Value = 3.483702037305
$ws.Cells.Item(4,19).Value = 0.0005879801490109799
$ws.Cells.Item(4,20).Value = 0.0005879801490109801
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Angptl3"
$ws.Cells.Item(5,3).Value = "Itga5"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.030331
$ws.Cells.Item(5,8).Value = 0.090993
$ws.Cells.Item(5,9).Value = 0.005348220947889333
$ws.Cells.Item(5,10).Value = 0.005348220947889334
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 16.966758
$ws.Cells.Item(5,14).Value = 50.900274
$ws.Cells.Item(5,15).Value = 0.1461640017030168
$ws.Cells.Item(5,16).Value = 0.1461640017030168
$ws.Cells.Item(5,17).Value = 0.514618736898
$ws.Cells.Item(5,18).Value = 4.631568632082
$ws.Cells.Item(5,19).Value = 0.0007817173757354066
$ws.Cells.Item(5,20).Value = 0.0007817173757354069
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Angptl3"
$ws.Cells.Item(6,3).Value = "Itga5"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3.081251
$ws.Cells.Item(6,8).Value = 9.243753
$ws.Cells.Item(6,9).Value = 0.5433124903202978
$ws.Cells.Item(6,10).Value = 0.5433124903202979
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 45.90594266666667
$ws.Cells.Item(6,14).Value = 137.717828
$ws.Cells.Item(6,15).Value = 0.3954672001633582
$ws.Cells.Item(6,16).Value = 0.3954672001633583
$ws.Cells.Item(6,17).Value = 141.4477317476093
$ws.Cells.Item(6,18).Value = 1273.029585728484
$ws.Cells.Item(6,19).Value = 0.2148622693607498
$ws.Cells.Item(6,20).Value = 0.2148622693607499
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Angptl3"
$ws.Cells.Item(7,3).Value = "Itga5"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 3.081251
$ws.Cells.Item(7,8).Value = 9.243753
$ws.Cells.Item(7,9).Value = 0.5433124903202978
$ws.Cells.Item(7,10).Value = 0.5433124903202979
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 40.44578266666667
$ws.Cells.Item(7,14).Value = 121.337348
$ws.Cells.Item(7,15).Value = 0.3484294080560655
$ws.Cells.Item(7,16).Value = 0.3484294080560656
$ws.Cells.Item(7,17).Value = 124.6236082874493
$ws.Cells.Item(7,18).Value = 1121.612474587044
$ws.Cells.Item(7,19).Value = 0.1893060493917682
$ws.Cells.Item(7,20).Value = 0.1893060493917683
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Angptl3"
$ws.Cells.Item(8,3).Value = "Itga5"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 3.081251
$ws.Cells.Item(8,8).Value = 9.243753
$ws.Cells.Item(8,9).Value = 0.5433124903202978
$ws.Cells.Item(8,10).Value = 0.5433124903202979
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 12.761795
$ws.Cells.Item(8,14).Value = 38.28538500000001
$ws.Cells.Item(8,15).Value = 0.1099393900775594
$ws.Cells.Item(8,16).Value = 0.1099393900775594
$ws.Cells.Item(8,17).Value = 39.32229360554501
$ws.Cells.Item(8,18).Value = 353.900642449905
$ws.Cells.Item(8,19).Value = 0.05973144380733344
$ws.Cells.Item(8,20).Value = 0.05973144380733347
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Angptl3"
$ws.Cells.Item(9,3).Value = "Itga5"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 3.081251
$ws.Cells.Item(9,8).Value = 9.243753
$ws.Cells.Item(9,9).Value = 0.5433124903202978
$ws.Cells.Item(9,10).Value = 0.5433124903202979
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 16.966758
$ws.Cells.Item(9,14).Value = 50.900274
$ws.Cells.Item(9,15).Value = 0.1461640017030168
$ws.Cells.Item(9,16).Value = 0.1461640017030168
$ws.Cells.Item(9,17).Value = 52.278840054258
$ws.Cells.Item(9,18).Value = 470.5095604883219
$ws.Cells.Item(9,19).Value = 0.07941272776044631
$ws.Cells.Item(9,20).Value = 0.07941272776044635
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Angptl3"
$ws.Cells.Item(10,3).Value = "Itga5"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.424064666666667
$ws.Cells.Item(10,8).Value = 4.272194
$ws.Cells.Item(10,9).Value = 0.2511032435928821
$ws.Cells.Item(10,10).Value = 0.2511032435928821
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 45.90594266666667
$ws.Cells.Item(10,14).Value = 137.717828
$ws.Cells.Item(10,15).Value = 0.3954672001633582
$ws.Cells.Item(10,16).Value = 0.3954672001633583
$ws.Cells.Item(10,17).Value = 65.37303094162577
$ws.Cells.Item(10,18).Value = 588.3572784746319
$ws.Cells.Item(10,19).Value = 0.09930309669561477
$ws.Cells.Item(10,20).Value = 0.09930309669561482
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Angptl3"
$ws.Cells.Item(11,3).Value = "Itga5"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.424064666666667
$ws.Cells.Item(11,8).Value = 4.272194
$ws.Cells.Item(11,9).Value = 0.2511032435928821
$ws.Cells.Item(11,10).Value = 0.2511032435928821
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 40.44578266666667
$ws.Cells.Item(11,14).Value = 121.337348
$ws.Cells.Item(11,15).Value = 0.3484294080560655
$ws.Cells.Item(11,16).Value = 0.3484294080560656
$ws.Cells.Item(11,17).Value = 57.5974100112791
$ws.Cells.Item(11,18).Value = 518.3766901015119
$ws.Cells.Item(11,19).Value = 0.08749175452602592
$ws.Cells.Item(11,20).Value = 0.08749175452602596
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Angptl3"
$ws.Cells.Item(12,3).Value = "Itga5"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.424064666666667
$ws.Cells.Item(12,8).Value = 4.272194
$ws.Cells.Item(12,9).Value = 0.2511032435928821
$ws.Cells.Item(12,10).Value = 0.2511032435928821
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 12.761795
$ws.Cells.Item(12,14).Value = 38.28538500000001
$ws.Cells.Item(12,15).Value = 0.1099393900775594
$ws.Cells.Item(12,16).Value = 0.1099393900775594
$ws.Cells.Item(12,17).Value = 18.17362134274333
$ws.Cells.Item(12,18).Value = 163.56259208469
$ws.Cells.Item(12,19).Value = 0.02760613744709828
$ws.Cells.Item(12,20).Value = 0.02760613744709829
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Angptl3"
$ws.Cells.Item(13,3).Value = "Itga5"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.424064666666667
$ws.Cells.Item(13,8).Value = 4.272194
$ws.Cells.Item(13,9).Value = 0.2511032435928821
$ws.Cells.Item(13,10).Value = 0.2511032435928821
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 16.966758
$ws.Cells.Item(13,14).Value = 50.900274
$ws.Cells.Item(13,15).Value = 0.1461640017030168
$ws.Cells.Item(13,16).Value = 0.1461640017030168
$ws.Cells.Item(13,17).Value = 24.161760575684
$ws.Cells.Item(13,18).Value = 217.455845181156
$ws.Cells.Item(13,19).Value = 0.03670225492414306
$ws.Cells.Item(13,20).Value = 0.03670225492414307
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Angptl3"
$ws.Cells.Item(14,3).Value = "Itga5"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 1.135585
$ws.Cells.Item(14,8).Value = 3.406755
$ws.Cells.Item(14,9).Value = 0.2002360451389308
$ws.Cells.Item(14,10).Value = 0.2002360451389308
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 45.90594266666667
$ws.Cells.Item(14,14).Value = 137.717828
$ws.Cells.Item(14,15).Value = 0.3954672001633582
$ws.Cells.Item(14,16).Value = 0.3954672001633583
$ws.Cells.Item(14,17).Value = 52.13009990312667
$ws.Cells.Item(14,18).Value = 469.17089912814
$ws.Cells.Item(14,19).Value = 0.07918678814287676
$ws.Cells.Item(14,20).Value = 0.07918678814287677
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Angptl3"
$ws.Cells.Item(15,3).Value = "Itga5"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 1.135585
$ws.Cells.Item(15,8).Value = 3.406755
$ws.Cells.Item(15,9).Value = 0.2002360451389308
$ws.Cells.Item(15,10).Value = 0.2002360451389308
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 40.44578266666667
$ws.Cells.Item(15,14).Value = 121.337348
$ws.Cells.Item(15,15).Value = 0.3484294080560655
$ws.Cells.Item(15,16).Value = 0.3484294080560656
$ws.Cells.Item(15,17).Value = 45.92962410952667
$ws.Cells.Item(15,18).Value = 413.36661698574
$ws.Cells.Item(15,19).Value = 0.06976812667924526
$ws.Cells.Item(15,20).Value = 0.06976812667924527
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Angptl3"
$ws.Cells.Item(16,3).Value = "Itga5"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 1.135585
$ws.Cells.Item(16,8).Value = 3.406755
$ws.Cells.Item(16,9).Value = 0.2002360451389308
$ws.Cells.Item(16,10).Value = 0.2002360451389308
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 12.761795
$ws.Cells.Item(16,14).Value = 38.28538500000001
$ws.Cells.Item(16,15).Value = 0.1099393900775594
$ws.Cells.Item(16,16).Value = 0.1099393900775594
$ws.Cells.Item(16,17).Value = 14.492102975075
$ws.Cells.Item(16,18).Value = 130.428926775675
$ws.Cells.Item(16,19).Value = 0.0220138286741167
$ws.Cells.Item(16,20).Value = 0.02201382867411671
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Angptl3"
$ws.Cells.Item(17,3).Value = "Itga5"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 1.135585
$ws.Cells.Item(17,8).Value = 3.406755
$ws.Cells.Item(17,9).Value = 0.2002360451389308
$ws.Cells.Item(17,10).Value = 0.2002360451389308
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 16.966758
$ws.Cells.Item(17,14).Value = 50.900274
$ws.Cells.Item(17,15).Value = 0.1461640017030168
$ws.Cells.Item(17,16).Value = 0.1461640017030168
$ws.Cells.Item(17,17).Value = 19.26719588343
$ws.Cells.Item(17,18).Value = 173.40476295087
$ws.Cells.Item(17,19).Value = 0.02926730164269202
$ws.Cells.Item(17,20).Value = 0.02926730164269203
